$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = 0.99834465265485828
$ws.Range("AW1").Value = 0.90005712158730766
$ws.Range("T2").Value = 0.77226798832662702
$ws.Range("AY2").Value = 0.87643816505791317
$ws.Range("AV4").Value = 0.67217950358074685
$ws.Range("I5").Value = 0.85291977505365746
$ws.Range("W6").Value = 0.83264815049043417
$ws.Range("AG6").Value = 0.88011583193298926
$ws.Range("F7").Value = 0.66165265461692924
$ws.Range("AC7").Value = 0.93800042302306275
$ws.Range("V8").Value = 0.70205545925806301
$ws.Range("BL9").Value = 0.88082768822846647
$ws.Range("BP9").Value = 0.86354063997574904
$ws.Range("R10").Value = 0.93529686474217688
$ws.Range("AY10").Value = 0.86988498074034681
$ws.Range("J11").Value = 0.59846397322674683
$ws.Range("P11").Value = 0.98281998579342567
$ws.Range("AT11").Value = 0.92187372399845902
$ws.Range("C12").Value = 0.67974851856674878
$ws.Range("BM12").Value = 0.87510284238655744
$ws.Range("BP12").Value = 0.95095331397937044
$ws.Range("E13").Value = 0.80178961575265428
$ws.Range("AB13").Value = 0.82850426714565095
$ws.Range("AO13").Value = 0.83758829612571495
$ws.Range("AX13").Value = 0.76492649182269834
$ws.Range("AA14").Value = 0.9196524602988565
$ws.Range("E15").Value = 0.90057425441872074
$ws.Range("Y15").Value = 0.80255103212412593
$ws.Range("D16").Value = 0.95789417627719597
$ws.Range("Z17").Value = 0.69563925291870154
$ws.Range("AU17").Value = 0.85995891487148701
$ws.Range("AM18").Value = 0.81094572714172042
$ws.Range("J19").Value = 0.99904683052185761
$ws.Range("V20").Value = 0.68476819397303679
$ws.Range("AV20").Value = 0.6995682478660904
$ws.Range("AK21").Value = 0.75834777055756564
$ws.Range("AN21").Value = 0.75034519476611639
$ws.Range("BM21").Value = 0.75089797590025031
$ws.Range("Q22").Value = 0.67991255652680915
$ws.Range("AM22").Value = 0.73064928024204856
$ws.Range("C23").Value = 0.93683644217955919
$ws.Range("U23").Value = 0.72699323638621616
$ws.Range("V23").Value = 0.75293692040184412
$ws.Range("AS23").Value = 0.98156474422176432
$ws.Range("AZ23").Value = 0.91823739188943732
$ws.Range("BO24").Value = 0.71901081471877071
$ws.Range("AA25").Value = 0.64215915313871208
$ws.Range("BE26").Value = 0.94276636734831187
$ws.Range("AB27").Value = 0.85339843207296551
$ws.Range("AC27").Value = 0.88339950615928231
$ws.Range("BO27").Value = 0.90368238276128721
$ws.Range("BC28").Value = 0.62501359970820913
$ws.Range("AQ29").Value = 0.75593277194420283
$ws.Range("M30").Value = 0.93430701165155783
$ws.Range("AC30").Value = 0.93772743337325459
$ws.Range("AK30").Value = 0.81242407624355928
$ws.Range("BN30").Value = 0.83266634176139886
$ws.Range("AC31").Value = 0.85268710019080884
$ws.Range("AH31").Value = 0.73795248211346764
$ws.Range("AY31").Value = 0.88896292029359802
$ws.Range("AA32").Value = 0.71715138045993965
$ws.Range("AO32").Value = 0.92534717977290337
$ws.Range("BB32").Value = 0.95054253104159936
$ws.Range("M33").Value = 0.86036301143492611
$ws.Range("AI33").Value = 0.76663787000646999
$ws.Range("BM33").Value = 0.84248699566262353
$ws.Range("A34").Value = 0.5290293059257023
$ws.Range("P34").Value = 0.94265738314335046
$ws.Range("AU34").Value = 0.67040440483189911
$ws.Range("BF34").Value = 0.77903820819939518
$ws.Range("AJ35").Value = 0.95797595293973126
$ws.Range("K36").Value = 0.96411032621373716
$ws.Range("M36").Value = 0.85247361465070004
$ws.Range("BE36").Value = 0.76875256896347899
$ws.Range("J37").Value = 0.63791138183788798
$ws.Range("G38").Value = 0.67546410147144731
$ws.Range("AW38").Value = 0.78408435902545492
$ws.Range("X39").Value = 0.98662042016881957
$ws.Range("BD40").Value = 0.77691897494905504
$ws.Range("BI40").Value = 0.97975715479563363
$ws.Range("AP41").Value = 0.78963581497153967
$ws.Range("F42").Value = 0.74790295329714418
$ws.Range("R42").Value = 0.80430716869264574
$ws.Range("BB42").Value = 0.85226500998609034
$ws.Range("AW43").Value = 0.73392735653890251
$ws.Range("AY43").Value = 0.78117691778350484
$ws.Range("Z44").Value = 0.76625809349849683
$ws.Range("AE44").Value = 0.88111633607467787
$ws.Range("AT44").Value = 0.92513519067067818
$ws.Range("BA44").Value = 0.84113093833308328
$ws.Range("O45").Value = 0.80763406493072087
$ws.Range("S45").Value = 0.82014408833864438
$ws.Range("AN46").Value = 0.90246214608998976
$ws.Range("J47").Value = 0.9526902700489378
$ws.Range("BC48").Value = 0.88954845600421728
$ws.Range("BP48").Value = 0.86381799560638228
$ws.Range("J49").Value = 0.7354844758081952
$ws.Range("AC49").Value = 0.85105058332402739
$ws.Range("E51").Value = 0.90028258763422664
$ws.Range("I51").Value = 0.87665867222574678
$ws.Range("Q51").Value = 0.87145542746889926
$ws.Range("L52").Value = 0.82266447750087723
$ws.Range("AY52").Value = 0.73357308873132332
$ws.Range("BA52").Value = 0.84077304778529394
$ws.Range("BF52").Value = 0.79928307041722613
$ws.Range("G53").Value = 0.7718156878138116
$ws.Range("AT53").Value = 0.89789944140969058
$ws.Range("L54").Value = 0.69108335782037633
$ws.Range("D55").Value = 0.8066995533635275
$ws.Range("BO55").Value = 0.88443834144347544
$ws.Range("F57").Value = 0.88015106340697491
$ws.Range("H57").Value = 0.69642787869131484
$ws.Range("AT57").Value = 0.84538314126991443
$ws.Range("BD58").Value = 0.51339478624678914
$ws.Range("AD59").Value = 0.82037752764993654
$ws.Range("R60").Value = 0.64825384277809306
$ws.Range("AI60").Value = 0.75563829627080414
$ws.Range("BK60").Value = 0.97977653564643563
$ws.Range("BG61").Value = 0.82317374530422427
$ws.Range("BH61").Value = 0.99312352057801601
$ws.Range("BP61").Value = 0.91852257240931334
$ws.Range("N62").Value = 0.58955989831860722
$ws.Range("Y62").Value = 0.93861856174351765
$ws.Range("BK62").Value = 0.90381172546066924
$ws.Range("AL63").Value = 0.679674009463066
$ws.Range("AU64").Value = 0.77534350920815531
$ws.Range("BK64").Value = 0.87528398107438066
$ws.Range("J66").Value = 0.93641406804660199
$ws.Range("BL66").Value = 0.99697488892257957
$ws.Range("BP66").Value = 0.92030682371297723
$ws.Range("S67").Value = 0.64375047770680949
$ws.Range("AF67").Value = 0.94738725105354815
$ws.Range("AL67").Value = 0.59919620593341349
$ws.Range("P68").Value = 0.74355790936812682
$ws.Range("AS68").Value = 0.52183301974987928
$ws.Range("AX68").Value = 0.54979866931841537
